# Fix: correção de erro ortográfico
#  - bump the cached "datetimeFigureOut" date placeholder text
#    (Slide Master + every Custom Layout) from 03/03/2021 -> 04/03/2021
#  - fix typo "Jorna de trabalho excessiva" -> "Jornada de trabalho excessiva"

$p = $ppt.ActivePresentation

$oldDate = "03/03/2021"
$newDate = "04/03/2021"

# ppPlaceholderDate = 16
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Type -eq 14 -and $shp.HasTextFrame -eq -1) {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# 1. Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 2. Every Custom Layout hanging off the master
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# 3. Typo fix on the actual slide content
for ($s = 1; $s -le $p.Slides.Count; $s++) {
    $slide = $p.Slides.Item($s)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "Jorna de trabalho excessiva") {
                $tr.Text = "Jornada de trabalho excessiva"
            }
        }
    }
}
